# Card17 — log a new maintenance event (row 16) and backfill the
# placeholder "nan" values that were left blank on row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# Row 15 previously had several untouched/blank columns — fill them with
# the same "nan" placeholder used elsewhere in this sheet.
$ws.Range("B15").Value = "nan"
$ws.Range("C15").Value = "nan"
$ws.Range("D15").Value = "nan"
$ws.Range("E15").Value = "nan"
$ws.Range("F15").Value = "nan"
$ws.Range("G15").Value = "nan"
$ws.Range("H15").Value = "nan"
$ws.Range("I15").Value = "nan"
$ws.Range("J15").Value = "nan"
$ws.Range("K15").Value = "nan"
$ws.Range("M15").Value = "nan"

# Append the new event as row 16 ("card" id stays text, matching the rest
# of column A).
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "17"
$ws.Range("A16").Style = "Normal"

$ws.Range("L16").Value = "27\4\2025"
$ws.Range("M16").Value = "593 t"
$ws.Range("N16").Value = "تم تغيير الجرائد الاماميه (1_2_4_5_7_8) ومعيارته"
$ws.Range("O16").Value = "الخبير"
